$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'42.905.90"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  -0.54%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'2.300.43"
$ws.Range("D3").Style = "Normal"
$ws.Range("E4").Value = "'  -0.04%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'299.74"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  -1.17%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'97.59"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  -2.06%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("E8").Value = "'  -0.02%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.506"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  -2.22%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'35.80"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  -1.24%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("E11").Value = "'  -0.98%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("E12").Value = "'  +0.66%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'17.70"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  -0.51%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("E14").Value = "'  -2.52%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'2.658.78"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  -0.99%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'2.291.05"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  -1.66%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("E17").Value = "'  -2.49%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'42.873.21"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  -0.45%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'12.61"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  -2.62%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("E20").Value = "'  -0.65%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("E21").Value = "'  -2.09%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'67.92"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  -0.44%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'241.90"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  +0.32%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("E24").Value = "'  -1.35%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("E25").Value = "'  +0.12%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'2.42"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  -1.89%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("E27").Value = "'  -0.30%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'25.01"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  -2.54%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'166.05"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  -1.56%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("E30").Value = "'  -0.56%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("E31").Value = "'  -1.87%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = "'32.62"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  -4.80%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("E33").Value = "'  -0.01%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = "'4.77"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  -3.54%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("E35").Value = "'  -3.07%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = "'17.43"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  -2.08%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("E37").Value = "'  -0.33%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = "'0.0686"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  -1.93%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'0.100"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  -2.15%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("E40").Value = "'  -3.12%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("E41").Value = "'  -1.14%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("E42").Value = "'  +0.32%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'2.000.15"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  +0.23%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("E44").Value = "'  -1.52%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("B45").Value = "'ApeXProtocol"
$ws.Range("B45").Style = "Normal"
$ws.Range("C45").Value = "'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("C45").Style = "Normal"
$ws.Range("D45").Value = "'2.16"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  -3.57%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("B46").Value = "'FraxShare"
$ws.Range("B46").Style = "Normal"
$ws.Range("C46").Value = "'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("C46").Style = "Normal"
$ws.Range("D46").Value = "'10.12"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  +0.14%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'17.26"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  -2.84%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D49").Value = "'2.525.70"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  -0.96%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'53.23"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  -3.38%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'72.11"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  -5.56%  "
$ws.Range("E51").Style = "Normal"
